$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (rows 18-21 = USDJPY block, rows 25-28 = XAUUSD block)
$data = @{
    18 = @{ B=65;  C=68;  D=71;  E=105; F=86;  G=116; H=511; I=32.67263427109975 }
    19 = @{ B=18;  C=32;  D=36;  E=37;  F=47;  G=87;  H=257; I=16.43222506393862 }
    20 = @{ B=123; C=101; D=122; E=108; F=95;  G=33;  H=582; I=37.21227621483376 }
    21 = @{ B=53;  C=58;  D=30;  E=12;  F=35;  G=26;  H=214; I=13.68286445012788 }
    25 = @{ B=72;  C=91;  D=85;  E=82;  F=64;  G=140; H=534; I=32.88177339901478 }
    26 = @{ B=21;  C=29;  D=56;  E=81;  F=42;  G=54;  H=283; I=17.42610837438423 }
    27 = @{ B=111; C=94;  D=86;  E=39;  F=111; G=104; H=545; I=33.55911330049261 }
    28 = @{ B=53;  C=44;  D=31;  E=57;  F=41;  G=36;  H=262; I=16.13300492610837 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}

$wb.Save()
